# ===== edit.ps1 : Weekly 112th Precinct CompStat update (crime data refresh) =====
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number and reporting week dates (rich-text shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "50"
$ws.Range("C9").Characters(27, 9).Text = "12/9/2024"
$ws.Range("C9").Characters(47, 9).Text = "12/15/2024"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 20
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 81
$ws.Range("J16").Value = 71
$ws.Range("K16").Value = 14.084507042253
$ws.Range("L16").Value = -7.954545454545
$ws.Range("M16").Value = -30.769230769230
$ws.Range("N16").Value = -87.04
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 109
$ws.Range("J17").Value = 97
$ws.Range("K17").Value = 12.371134020618
$ws.Range("L17").Value = 34.567901234567
$ws.Range("M17").Value = 91.228070175438
$ws.Range("N17").Value = -11.382113821138
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 98
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -14.035087719298
$ws.Range("L18").Value = -15.517241379310
$ws.Range("M18").Value = -22.834645669291
$ws.Range("N18").Value = -92.648162040510
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 37.5
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -37.777777777777
$ws.Range("I19").Value = 432
$ws.Range("J19").Value = 465
$ws.Range("K19").Value = -7.096774193548
$ws.Range("L19").Value = -10.187110187110
$ws.Range("M19").Value = 16.756756756756
$ws.Range("N19").Value = -55.601233299075
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -45.454545454545
$ws.Range("I20").Value = 159
$ws.Range("J20").Value = 163
$ws.Range("K20").Value = -2.453987730061
$ws.Range("L20").Value = 55.882352941176
$ws.Range("M20").Value = 63.917525773195
$ws.Range("N20").Value = -95.097132284921
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 31.25
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 891
$ws.Range("J21").Value = 916
$ws.Range("K21").Value = -2.729257641921
$ws.Range("L21").Value = 0.791855203619
$ws.Range("M21").Value = 15.265200517464
$ws.Range("N21").Value = -85.884030418250
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 37
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = 32.142857142857
$ws.Range("L22").Value = 2.777777777777
$ws.Range("M22").Value = 68.181818181818
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -2.857142857142
$ws.Range("F24").Value = 166
$ws.Range("H24").Value = 44.347826086956
$ws.Range("I24").Value = 1708
$ws.Range("J24").Value = 1483
$ws.Range("K24").Value = 15.171948752528
$ws.Range("L24").Value = -0.870574579222
$ws.Range("M24").Value = 80.740740740740
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -4
$ws.Range("F25").Value = 123
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = 51.851851851851
$ws.Range("I25").Value = 1259
$ws.Range("J25").Value = 1057
$ws.Range("K25").Value = 19.110690633869
$ws.Range("L25").Value = 3.112203112203
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 166.666666666667
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 57.142857142857
$ws.Range("I26").Value = 267
$ws.Range("J26").Value = 233
$ws.Range("K26").Value = 14.592274678111
$ws.Range("L26").Value = 45.901639344262
$ws.Range("M26").Value = 27.751196172248
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("L27").Value = 5
$ws.Range("F28").Value = 4
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 8.333333333333
$ws.Range("L28").Value = -2.5
$ws.Range("I31").Value = 11
$ws.Range("K31").Value = -21.428571428571
$ws.Range("L31").Value = 37.5

# --- Cells that change from placeholder text ("0"/"***.*") to real numeric values ---
# (style must be converted from General/text (13) to the numeric format used by the rest of the column)
function Set-NumericCell($addr, $val, $formatSource) {
    $ws.Range($addr).Value = $val
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}

Set-NumericCell "C20" 3 "D17"
Set-NumericCell "D22" 2 "D17"
Set-NumericCell "E22" -50 "E17"
Set-NumericCell "G22" 2 "D17"
Set-NumericCell "H22" 100 "E17"
Set-NumericCell "D28" 1 "D17"
Set-NumericCell "E28" -100 "E17"
Set-NumericCell "G28" 1 "D17"
Set-NumericCell "H28" 300 "E17"

# --- Cell that changes from a numeric value back to the placeholder text "0" ---
function Set-PlaceholderTextCell($addr, $text, $formatSource) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
}
Set-PlaceholderTextCell "C28" "0" "D27"
